$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append ".60270.60134" to the existing customer code list stored in B6
$ws.Range("B6").Value = $ws.Range("B6").Value2 + ".60270.60134"

# Update the sheet view (scroll position and active selection)
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("B7").Select()
